# Updated cryptos list - applies the per-row coin/price/volume changes
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure Price (D) and Volume(1h) (E) columns stay text, exactly like the
# source data (inline strings), so values such as "12.00" or "28.176.21"
# are not silently reinterpreted as numbers by Excel.
$ws.Range("D2:E51").NumberFormat = "@"

$updates = @(
    @{Row=2;  D='28.176.21';   E='  +0.81%  '}
    @{Row=3;  D='1.801.25';    E='  +2.40%  '}
    @{Row=4;  D='1.003';       E='  -0.02%  '}
    @{Row=5;  D='336.67';      E='  +0.17%  '}
    @{Row=6;  D='0.9996';      E='  +0.03%  '}
    @{Row=7;  D='0.4617';      E='  +20.57%  '}
    @{Row=8;  D='0.3699';      E='  +9.21%  '}
    @{Row=9;  D='45.11';       E='  +0.61%  '}
    @{Row=10; D='1.147';       E='  +3.04%  '}
    @{Row=11; D='0.07577';     E='  +4.98%  '}
    @{Row=12;                  E='  +0.00%  '}
    @{Row=13; D='22.29';       E='  -0.06%  '}
    @{Row=14; D='6.326';       E='  +3.11%  '}
    @{Row=15; D='7.421';       E='  +3.54%  '}
    @{Row=16; D='1.800.35';    E='  +2.43%  '}
    @{Row=17; D='0.00001095';  E='  +3.63%  '}
    @{Row=18; D='0.06716';     E='  +1.59%  '}
    @{Row=19; D='82.03';       E='  +3.56%  '}
    @{Row=20;                  E='  -0.04%  '}
    @{Row=21; D='17.48';       E='  +5.12%  '}
    @{Row=22; D='6.395';       E='  +2.76%  '}
    @{Row=23; D='28.174.46';   E='  +0.75%  '}
    @{Row=24; D='11.85';       E='  +2.00%  '}
    @{Row=25; D='2.417';       E='  +1.04%  '}
    @{Row=26; D='20.66';       E='  +4.57%  '}
    @{Row=27; B='Monero';        C='https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr';           D='152.79';    E='  +0.38%  '}
    @{Row=28; B='LidoDAOToken';  C='https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo';          D='2.373';     E='  +2.84%  '}
    @{Row=29; D='2.005.66';    E='  +2.36%  '}
    @{Row=30; D='133.50';      E='  +1.44%  '}
    @{Row=31; D='1.251';       E='  -1.31%  '}
    @{Row=32; D='4.030';       E='  +0.40%  '}
    @{Row=33; D='0.09569';     E='  +8.52%  '}
    @{Row=34; D='5.846';       E='  +0.22%  '}
    @{Row=35; D='0.2219';      E='  +5.19%  '}
    @{Row=36; D='0.06354';     E='  +2.17%  '}
    @{Row=37; D='0.02345';     E='  +2.43%  '}
    @{Row=38; D='12.00';       E='  -1.56%  '}
    @{Row=39; D='5.238';       E='  +1.40%  '}
    @{Row=40; D='0.6621';      E='  +0.19%  '}
    @{Row=41; D='1.507';       E='  +1.92%  '}
    @{Row=42; D='1.232';       E='  +1.79%  '}
    @{Row=43; D='8.061';       E='  +1.25%  '}
    @{Row=44; B='EnergySwap';  C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens';             D='14.15';     E='  +2.76%  '}
    @{Row=45; B='Frax';        C='https://coinranking.com/coin/KfWtaeV1W+frax-frax';                  D='0.9996';    E='  -0.02%  '}
    @{Row=46; D='0.6088';      E='  +0.95%  '}
    @{Row=47; D='3.833';       E='  +0.38%  '}
    @{Row=48; D='129.88';      E='  +2.92%  '}
    @{Row=49; D='2.042';       E='  +1.72%  '}
    @{Row=50; B='Cronos';      C='https://coinranking.com/coin/65PHZTpmE55b+cronos-cro';               D='0.07145';   E='  +2.18%  '}
    @{Row=51; B='EOS';         C='https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos';                 D='1.173';     E='  -0.04%  '}
)

foreach ($u in $updates) {
    $r = $u.Row
    if ($u.ContainsKey('B')) { $ws.Range("B$r").Value = $u.B }
    if ($u.ContainsKey('C')) { $ws.Range("C$r").Value = $u.C }
    if ($u.ContainsKey('D')) { $ws.Range("D$r").Value = $u.D }
    if ($u.ContainsKey('E')) { $ws.Range("E$r").Value = $u.E }
}
